$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$names = @(
    'NU Hospitals',
    'Jayashree Multispeciality Hospital',
    'Springleaf Healthcare',
    'NU Hospitals',
    'Trustwell Hospitals',
    'HCG - The Specialist in Cancer Care',
    'ReSCUE Hospital',
    'Navachethana Hospital',
    'Motherhood Hospital',
    'Narayana Hrudhayalaya Multispeciality Hospital',
    'Koshys Hospital',
    'Sakra World Hospital',
    'Poornima Hospital',
    'Cloudnine Children''s Hospital - Jayanagar',
    'Aaxis Super Speciality Hospital',
    'Specialist Hospital',
    'Columbia Asia Hospital',
    'Greenview Medical Center',
    'St. Philomena''s Hospital',
    'Axon Speciality Hospital',
    'Blossom multispecialty hospital',
    'Sagar Hospitals',
    'Ramakrishna Super Speciality Hospital',
    'Cloudnine Hospital - Old Airport Road',
    'Motherhood Hospital',
    'Queens Hospital',
    'Spandana Hospital & HCC',
    'Sanjeevini MultispecialIty Hospital',
    'Mallige Medical Centre',
    'Divakars Speciality Hospital',
    'Cloudnine Hospital - Whitefield',
    'Manipal Hospital',
    'RxDx Healthcare',
    'Birthright By Rainbow',
    'Cloudnine Hospital - Bellandur',
    'Cloudnine Hospital - Jayanagar',
    'Apple Hospital',
    'Ovum Woman & Child Speciality Hospital',
    'Aveksha Hospital',
    'Brindhavvan Areion Hospital',
    'Sakhii Hospital',
    'Motherhood Hospital',
    'Cloudnine Fertility - IVF Centre, Old Airport Road',
    'Devagiri Hospital',
    'Republic Hospital - Powered by Medisync',
    'SPARSH Hospital',
    'Acura Speciality Hospital',
    'Cloudnine Fertility - IVF Centre, Jayanagar',
    'Cloudnine Fertility - IVF Centre, Whitefield',
    'Bhagwan Mahaveer Jain Hospital',
    'Motherhood Hospital',
    'Spandana Hospital Rehabilitation Center',
    'Parkmed',
    'Birthright By Rainbow',
    'Punya Hospital',
    'Cloudnine Hospital - Malleshwaram',
    'Nethradhama Super Speciality Eye Hospital',
    'Wellspring Mother & Child Hospital',
    'AyurVAID Hospitals',
    'Abhaya Hospital',
    'Cloudnine Fertility - IVF Centre, Malleshwaram',
    'Rainbow Children''s Hospital',
    'S J Hospital',
    'Mazumdar Shaw Medical Center',
    'Family Hospital',
    'Ovum Woman & Child Speciality Hospital',
    'Cloudnine Fertility - IVF Centre, Bellandur',
    'Shekhar Hospital',
    'Acharya Tulasi Jain Hospital',
    'P D Hinduja Sindhi Hospital',
    'Cloudnine Fertility - IVF Centre, HRBR Layout',
    'Motherhood Fertility - Sarjapur Road',
    'Gunasheela Surgical and Maternity Hospital',
    'Prakriya Hospitals',
    'Al Ameen Hospital',
    'Shree Ayur Ganga',
    'Ovum Birthing Center',
    'Ovum Birthing Center',
    'Spandana Hospital Rehabilitation Center',
    'Comfort Multispeciality Hospital',
    'Parkmed',
    'Suraksha Multi Speciality Hospital',
    'Promed Hospital',
    'Birthright By Rainbow',
    'Seventh Day Adventist Hospital',
    'Punya Hospital',
    'Cloudnine Hospital - Malleshwaram',
    'Nethradhama Super Speciality Eye Hospital',
    'Wellspring Mother & Child Hospital',
    'AyurVAID Hospitals',
    'Abhaya Hospital',
    'Cloudnine Fertility - IVF Centre, Malleshwaram',
    'Sri Sai Nothside Hospital',
    'Sri Sai Krupa Hospital',
    'Deeksha Hospital',
    'Rainbow Children''s Hospital',
    'S J Hospital',
    'Mazumdar Shaw Medical Center',
    'Family Hospital',
    'Ovum Woman & Child Speciality Hospital',
    'Maharaja Agrasen Hospital',
    'Cloudnine Fertility - IVF Centre, Bellandur',
    'Rajashri Grandhim Speciality Hospital',
    'Vijayashree Multi Specialty Hospital',
    'Asha Nursing Home',
    'Shekhar Hospital',
    'Sevakshetra Hospital',
    'Acharya Tulasi Jain Hospital',
    'P D Hinduja Sindhi Hospital',
    'Cloudnine Fertility - IVF Centre, HRBR Layout',
    'Geetha Hospital',
    'Rajalakshmi Multispeciality Hospital.',
    'Sri Devi Hospital',
    'Motherhood Fertility - Sarjapur Road',
    'Avhieta Health Care and Hospital',
    'Gunasheela Surgical and Maternity Hospital',
    'Prakriya Hospitals',
    'Al Ameen Hospital',
    'Ovum Birthing Center',
    'First Eat Right Clinic',
    'Kanva Diagnostic Services Pvt Ltd',
    'First Eat Right Clinic',
    'Apollo Medical Centre',
    'Shree Ayur Ganga',
    'First Eat Right Clinic- Jayanagar 4 Block',
    'First Eat Right Clinic - Koramangala'
)

for ($i = 0; $i -lt $names.Length; $i++) {
    $ws.Cells.Item($i + 1, 1).Value = $names[$i]
}
